$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the urbanization labels so they distinguish the settlement-type
# category ("City"/"Village") from the broader urban/rural area labels.
$ws.Range("A23").Value = "Шаар жерлери"
$ws.Range("B23").Value = "Городские поселения"
$ws.Range("C23").Value = "City"

$ws.Range("A24").Value = "Айыл аймагы"
$ws.Range("B24").Value = "Сельская местность"
$ws.Range("C24").Value = "Village"

# Move the active selection to match the published workbook.
$null = $ws.Range("C30").Select()
